$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = -6
    4  = 0
    5  = -2
    8  = -4
    9  = 1
    10 = -1
    11 = -2
    12 = -4
    13 = -3
    14 = -4
    17 = 4
    18 = 3
    19 = -3
    20 = 2
    21 = -1
    22 = 7
    23 = 6
    24 = 4
    25 = 4
    26 = 5
    27 = 1
    28 = 2
    29 = -2
    31 = 4
    32 = 5
    33 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
